$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 10.79022983580887
$ws.Range("D2").Value = 9.475165924678562
$ws.Range("E2").Value = 13.91678701490353
$ws.Range("F2").Value = 32.09718107448681
$ws.Range("G2").Value = 32.96140971147478
$ws.Range("H2").Value = 15.52060635840146
$ws.Range("I2").Value = 25.17226234293376
$ws.Range("J2").Value = 10.16169778273283
$ws.Range("K2").Value = 17.73659794296289
$ws.Range("L2").Value = 10.37018813975018
$ws.Range("O2").Value = 24.1476135503099
$ws.Range("C3").Value = 10.72102891871319
$ws.Range("D3").Value = 9.420236036357634
$ws.Range("E3").Value = 13.90291528193176
$ws.Range("F3").Value = 32.25345759410513
$ws.Range("G3").Value = 33.20925651238893
$ws.Range("H3").Value = 15.59771430635909
$ws.Range("I3").Value = 25.29730471369855
$ws.Range("J3").Value = 10.1812461113481
$ws.Range("K3").Value = 17.04606431142616
$ws.Range("L3").Value = 10.38407781869823
$ws.Range("O3").Value = 24.29423069403494
$ws.Range("C4").Value = 10.67979973494859
$ws.Range("D4").Value = 9.387320267379634
$ws.Range("E4").Value = 13.8962444352582
$ws.Range("F4").Value = 32.35780198137804
$ws.Range("G4").Value = 33.37312099131682
$ws.Range("H4").Value = 15.64789720532501
$ws.Range("I4").Value = 25.3800098438682
$ws.Range("J4").Value = 10.19426596534733
$ws.Range("K4").Value = 16.60646229718066
$ws.Range("L4").Value = 10.39362290312193
$ws.Range("O4").Value = 24.39010207860801
$ws.Range("C5").Value = 10.66332860396209
$ws.Range("D5").Value = 9.374120115396359
$ws.Range("E5").Value = 13.89399296503783
$ws.Range("F5").Value = 32.40242715969512
$ws.Range("G5").Value = 33.44282285128294
$ws.Range("H5").Value = 15.66906141213004
$ws.Range("I5").Value = 25.41520114423518
$ws.Range("J5").Value = 10.19982781056113
$ws.Range("K5").Value = 16.42360797839796
$ws.Range("L5").Value = 10.39776871663675
$ws.Range("O5").Value = 24.43063970946301
$ws.Range("C6").Value = 10.66061389681895
$ws.Range("D6").Value = 9.371941387373019
$ws.Range("E6").Value = 13.89364738612111
$ws.Range("F6").Value = 32.4099640444434
$ws.Range("G6").Value = 33.4545731344628
$ws.Range("H6").Value = 15.67261886933885
$ws.Range("I6").Value = 25.42113444905747
$ws.Range("J6").Value = 10.20076683291621
$ws.Range("K6").Value = 16.39302763697912
$ws.Range("L6").Value = 10.3984726071482
$ws.Range("O6").Value = 24.43745965780722
$ws.Range("C7").Value = 10.67957624595469
$ws.Range("D7").Value = 9.387141369882029
$ws.Range("E7").Value = 13.89621217712147
$ws.Range("F7").Value = 32.35839530076225
$ws.Range("G7").Value = 33.37404918531093
$ws.Range("H7").Value = 15.64817974043057
$ws.Range("I7").Value = 25.38047842333357
$ws.Range("J7").Value = 10.19433993668181
$ws.Range("K7").Value = 16.60401099316211
$ws.Range("L7").Value = 10.3936777774404
$ws.Range("O7").Value = 24.39064283586061
$ws.Range("C8").Value = 10.76611577647616
$ws.Range("D8").Value = 9.456063326644939
$ws.Range("E8").Value = 13.91162245937945
$ws.Range("F8").Value = 32.14932044658381
$ws.Range("G8").Value = 33.04443427875535
$ws.Range("H8").Value = 15.54660441855545
$ws.Range("I8").Value = 25.21414495907432
$ws.Range("J8").Value = 10.16822718260138
$ws.Range("K8").Value = 17.50186181050254
$ws.Range("L8").Value = 10.3747665384468
$ws.Range("O8").Value = 24.19695266271707
$ws.Range("C9").Value = 10.94522328874102
$ws.Range("D9").Value = 9.59723959906405
$ws.Range("E9").Value = 13.95636892771991
$ws.Range("F9").Value = 31.80615035206049
$ws.Range("G9").Value = 32.4913544692944
$ws.Range("H9").Value = 15.36991082199041
$ws.Range("I9").Value = 24.93511864814448
$ws.Range("J9").Value = 10.12507305883669
$ws.Range("K9").Value = 19.13059847827982
$ws.Range("L9").Value = 10.34572979682794
$ws.Range("O9").Value = 23.86358294645079
$ws.Range("C10").Value = 11.08172103612097
$ws.Range("D10").Value = 9.704058855606691
$ws.Range("E10").Value = 13.99792321145488
$ws.Range("F10").Value = 31.59510668462714
$ws.Range("G10").Value = 32.1426801226538
$ws.Range("H10").Value = 15.25377208967715
$ws.Range("I10").Value = 24.75901768761951
$ws.Range("J10").Value = 10.09825342334112
$ws.Range("K10").Value = 20.23759882621242
$ws.Range("L10").Value = 10.32927657319878
$ws.Range("O10").Value = 23.64705547213559
$ws.Range("C11").Value = 11.14470094527899
$ws.Range("D11").Value = 9.753200338005584
$ws.Range("E11").Value = 14.01866978443611
$ws.Range("F11").Value = 31.50809154087982
$ws.Range("G11").Value = 31.99675922229986
$ws.Range("H11").Value = 15.20390065264063
$ws.Range("I11").Value = 24.6852138230532
$ws.Range("J11").Value = 10.08710832521912
$ws.Range("K11").Value = 20.72020805278015
$ws.Range("L11").Value = 10.32284542433961
$ws.Range("O11").Value = 23.55473710432893
$ws.Range("C12").Value = 11.16866158185828
$ws.Range("D12").Value = 9.771877016821144
$ws.Range("E12").Value = 14.02678707626977
$ws.Range("F12").Value = 31.47643993664601
$ws.Range("G12").Value = 31.94334321733729
$ws.Range("H12").Value = 15.18544099206389
$ws.Range("I12").Value = 24.65817588227188
$ws.Range("J12").Value = 10.08303931287767
$ws.Range("K12").Value = 20.89983982277631
$ws.Range("L12").Value = 10.32056110014056
$ws.Range("O12").Value = 23.5206693275509
$ws.Range("C13").Value = 11.16349652582794
$ws.Range("D13").Value = 9.767851818478688
$ws.Range("E13").Value = 14.02502733033499
$ws.Range("F13").Value = 31.48319879452858
$ws.Range("G13").Value = 31.95476517219156
$ws.Range("H13").Value = 15.18939768517069
$ws.Range("I13").Value = 24.66395846243534
$ws.Range("J13").Value = 10.08390891986035
$ws.Range("K13").Value = 20.86129328163104
$ws.Range("L13").Value = 10.32104636172888
$ws.Range("O13").Value = 23.5279667583506
$ws.Range("C14").Value = 11.14667005821565
$ws.Range("D14").Value = 9.754735583119063
$ws.Range("E14").Value = 14.01933239176631
$ws.Range("F14").Value = 31.50546147099665
$ws.Range("G14").Value = 31.99232766580009
$ws.Range("H14").Value = 15.20237343702057
$ws.Range("I14").Value = 24.68297113913829
$ws.Range("J14").Value = 10.08677053299652
$ws.Range("K14").Value = 20.73504958736762
$ws.Range("L14").Value = 10.32265446790719
$ws.Range("O14").Value = 23.55191645323418
$ws.Range("C15").Value = 11.13637740749889
$ws.Range("D15").Value = 9.746710018711022
$ws.Range("E15").Value = 14.01587794536799
$ws.Range("F15").Value = 31.51926737961544
$ws.Range("G15").Value = 32.01557601498777
$ws.Range("H15").Value = 15.21037687716954
$ws.Range("I15").Value = 24.69473555677745
$ws.Range("J15").Value = 10.08854305911915
$ws.Range("K15").Value = 20.65731223043938
$ws.Range("L15").Value = 10.32365913111662
$ws.Range("O15").Value = 23.56670246346415
$ws.Range("C16").Value = 11.07762162712202
$ws.Range("D16").Value = 9.700857445316757
$ws.Range("E16").Value = 13.9966041249652
$ws.Range("F16").Value = 31.60097477616216
$ws.Range("G16").Value = 32.152473158478
$ws.Range("H16").Value = 15.25709086605676
$ws.Range("I16").Value = 24.76396810318078
$ws.Range("J16").Value = 10.09900298655284
$ws.Range("K16").Value = 20.2056277458526
$ws.Range("L16").Value = 10.32971802239238
$ws.Range("O16").Value = 23.6532132721437
$ws.Range("C17").Value = 11.04179258359434
$ws.Range("D17").Value = 9.672861114802096
$ws.Range("E17").Value = 13.98524956456917
$ws.Range("F17").Value = 31.65340669260443
$ws.Range("G17").Value = 32.23971725839493
$ws.Range("H17").Value = 15.28650654803486
$ws.Range("I17").Value = 24.80805730336892
$ws.Range("J17").Value = 10.10568984637757
$ws.Range("K17").Value = 19.92308212679968
$ws.Range("L17").Value = 10.33370444680658
$ws.Range("O17").Value = 23.70786960968346
$ws.Range("C18").Value = 11.02126906282886
$ws.Range("D18").Value = 9.656810696596935
$ws.Range("E18").Value = 13.97889248017512
$ws.Range("F18").Value = 31.6844099145066
$ws.Range("G18").Value = 32.29109149899683
$ws.Range("H18").Value = 15.30370434876479
$ws.Range("I18").Value = 24.83400953621634
$ws.Range("J18").Value = 10.10963530142125
$ws.Range("K18").Value = 19.75860077640657
$ws.Range("L18").Value = 10.33609652325974
$ws.Range("O18").Value = 23.73988813242613
$ws.Range("C19").Value = 11.01433513039503
$ws.Range("D19").Value = 9.651385624118658
$ws.Range("E19").Value = 13.97677004373626
$ws.Range("F19").Value = 31.6950521373258
$ws.Range("G19").Value = 32.30869048143367
$ws.Range("H19").Value = 15.30957509765971
$ws.Range("I19").Value = 24.84289831516871
$ws.Range("J19").Value = 10.11098823895517
$ws.Range("K19").Value = 19.70257561002048
$ws.Range("L19").Value = 10.33692348971276
$ws.Range("O19").Value = 23.75082888198175
$ws.Range("C20").Value = 11.04559802675479
$ws.Range("D20").Value = 9.675836035803329
$ws.Range("E20").Value = 13.98644032223808
$ws.Range("F20").Value = 31.64773764465727
$ws.Range("G20").Value = 32.23030631311053
$ws.Range("H20").Value = 15.28334635762499
$ws.Range("I20").Value = 24.80330250566031
$ws.Range("J20").Value = 10.10496773892924
$ws.Range("K20").Value = 19.95336415094729
$ws.Range("L20").Value = 10.33326982270419
$ws.Range("O20").Value = 23.70199113466667
$ws.Range("C21").Value = 11.15160950097414
$ws.Range("D21").Value = 9.758586386074267
$ws.Range("E21").Value = 14.0209980827605
$ws.Range("F21").Value = 31.49888707397995
$ws.Range("G21").Value = 31.98124455700441
$ws.Range("H21").Value = 15.19855059545968
$ws.Range("I21").Value = 24.67736193222358
$ws.Range("J21").Value = 10.08592590189138
$ws.Range("K21").Value = 20.77221591759699
$ws.Range("L21").Value = 10.32217803361528
$ws.Range("O21").Value = 23.54485764142028
$ws.Range("C22").Value = 11.22153663071869
$ws.Range("D22").Value = 9.813058918237433
$ws.Range("E22").Value = 14.0451029417733
$ws.Range("F22").Value = 31.40917994853201
$ws.Range("G22").Value = 31.82920702505237
$ws.Range("H22").Value = 15.14561216167341
$ws.Range("I22").Value = 24.60035809255075
$ws.Range("J22").Value = 10.07436323667738
$ws.Range("K22").Value = 21.28914749729698
$ws.Range("L22").Value = 10.31580888982391
$ws.Range("O22").Value = 23.44735750816413
$ws.Range("C23").Value = 11.18416171525633
$ws.Range("D23").Value = 9.783953854725759
$ws.Range("E23").Value = 14.03210009674175
$ws.Range("F23").Value = 31.45636307959783
$ws.Range("G23").Value = 31.90936441618544
$ws.Range("H23").Value = 15.17363947988388
$ws.Range("I23").Value = 24.64096993741494
$ws.Range("J23").Value = 10.08045383831988
$ws.Range("K23").Value = 21.01495111731395
$ws.Range("L23").Value = 10.31912786602087
$ws.Range("O23").Value = 23.49891894958506
$ws.Range("C24").Value = 11.04387735086044
$ws.Range("D24").Value = 9.674490932882881
$ws.Range("E24").Value = 13.98590144836653
$ws.Range("F24").Value = 31.65029794468612
$ws.Range("G24").Value = 32.23455721348927
$ws.Range("H24").Value = 15.2847741877698
$ws.Range("I24").Value = 24.80545026613469
$ws.Range("J24").Value = 10.10529388879245
$ws.Range("K24").Value = 19.93968000033817
$ws.Range("L24").Value = 10.33346600402602
$ws.Range("O24").Value = 23.70464693728614
$ws.Range("C25").Value = 10.89584707337991
$ws.Range("D25").Value = 9.558461490667945
$ws.Range("E25").Value = 13.94272599049386
$ws.Range("F25").Value = 31.89179708090333
$ws.Range("G25").Value = 32.63090798801179
$ws.Range("H25").Value = 15.41530666854353
$ws.Range("I25").Value = 25.0055388317799
$ws.Range("J25").Value = 10.13588770147167
$ws.Range("K25").Value = 18.70514017759158
$ws.Range("L25").Value = 10.35272605880337
$ws.Range("O25").Value = 23.9487871581479
